# Add stock info rows for: TXN, ULVR.L, VZ, HEN3.DE, MPW, BA, HL.L, TIP,
# MCHP, ITW, TLT, KR, AMAT, XLV, NXPI, BPY, GRMN, ATVI, CME, AGN
# Columns are Symbol | ISIN | Name | Address | CountryCode (matches Table_1).

$wb  = $excel.ActiveWorkbook
$ws  = $wb.ActiveSheet
$tbl = $ws.ListObjects.Item(1)

$newRecords = @(
    "TXN|US8825081040|Texas Instruments Incorporated|12500 TI Boulevard Dallas, TX 75243 United States|US",
    "ULVR.L|GB00B10RZP78|Unilever PLC|Unilever House 100 Victoria Embankment London EC4Y 0DY United Kingdom|GB",
    "VZ|US92343V1044|Verizon Communications Inc.|1095 Avenue of the Americas, New York, NY 10036|US",
    "HEN3.DE|DE0006048408|Henkel AG & Co. KGaA|Henkelstrasse 67, Düsseldorf 40589, Germany|DE",
    "MPW|US58463J3041|Medical Properties Trust, Inc.|1000 Urban Center Drive, Suite 501, Birmingham, AL 35242, United States|US",
    "BA|US0970231058|The Boeing Company|100 North Riverside Plaza, Chicago, IL 60606-1596, United States|US",
    "HL.L|GB00B1VZ0M25|Hargreaves Lansdown plc|One College Square South, Anchor Road, Bristol BS1 5HL, United Kingdom|GB",
    "TIP|US4642871762|iShares TIPS Bond ETF|iShares|US",
    "MCHP|US5950171042|Microchip Technology Incorporated|2355 West Chandler Boulevard, Chandler, AZ 85224-6199, United States|US",
    "ITW|US4523081093|Illinois Tool Works Inc.|155 Harlem Avenue, Glenview, IL 60025-4075, United States|US",
    "TLT|US4642874329|iShares 20+ Year Treasury Bond ETF|iShares|US",
    "KR|US501044DN88|The Kroger Co.|1014 Vine Street, Cincinnati, OH 45202, United States|US",
    "AMAT|US0382221051|Applied Materials, Inc.|3050 Bowers Avenue, PO Box 58039, Santa Clara, CA 95052-3299, United States|US",
    "XLV|US81369Y2090|Health Care Select Sector SPDR Fund|SPDR State Street Global Advisors|US",
    "NXPI|NL0009538784|NXP Semiconductors N.V.|High Tech Campus 60, Eindhoven 5656 AG, Netherlands|NL",
    "BPY|BMG162491077|Brookfield Property Partners L.P.|73 Front Street, 5th Floor, Hamilton HM 12, Bermuda|BM",
    "GRMN|CH0114405324|Garmin Ltd.|Mühlentalstrasse 2, Schaffhausen 8200, Switzerland|CH",
    "ATVI|US00507V1098|Activision Blizzard, Inc.|3100 Ocean Park Boulevard, Santa Monica, CA 90405, United States|US",
    "CME|US12572Q1058|CME Group Inc.|20 South Wacker Drive, Chicago, IL 60606, United States|US",
    "AGN|CA01559R1038|Algernon Pharmaceuticals Inc.|700 West Pender Street, Suite 915, Vancouver, BC V6C 1G8, Canada|CA"
)

foreach ($record in $newRecords) {
    $fields = $record.Split("|")

    $newRow = $tbl.ListRows.Add()
    $r = $newRow.Range.Row

    $ws.Cells.Item($r, 1).Value = $fields[0]
    $ws.Cells.Item($r, 2).Value = $fields[1]
    $ws.Cells.Item($r, 3).Value = $fields[2]
    $ws.Cells.Item($r, 4).Value = $fields[3]

    $countryCell = $ws.Cells.Item($r, 5)
    $countryCell.Value = $fields[4]
    $countryCell.HorizontalAlignment = -4108
}
